# Update cryptos list: refresh Price (column D) and Volume(1h) (column E) values
# per the Jan 20 2024 GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.523.44'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '2.468.29'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.65'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.85'
$ws.Range("E6").Value = '  -3.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.549'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.511'
$ws.Range("E9").Value = '  +2.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.21'
$ws.Range("E10").Value = '  -4.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").Value = '2.847.46'
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.83'
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("E15").Value = '  +2.93%  '
$ws.Range("D16").Value = '2.473.69'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.771'
$ws.Range("E17").Value = '  -3.16%  '
$ws.Range("D18").Value = '41.510.06'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.47'
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("D20").Value = '0.0₃0941'
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.22'
$ws.Range("E21").Value = '  +3.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.06'
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.62'
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.90'
$ws.Range("E26").Value = '  -0.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.62'
$ws.Range("E27").Value = '  +1.52%  '
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.67'
$ws.Range("E29").Value = '  -1.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.35'
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.60'
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.43'
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.30'
$ws.Range("E35").Value = '  -4.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.87'
$ws.Range("E36").Value = '  -7.40%  '
$ws.Range("E37").Value = '  +1.42%  '
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.77'
$ws.Range("E39").Value = '  -6.20%  '
$ws.Range("E40").Value = '  -11.62%  '
$ws.Range("E41").Value = '  -4.48%  '
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("D43").Value = '1.940.74'
$ws.Range("E43").Value = '  -2.99%  '
$ws.Range("E44").Value = '  -1.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.33'
$ws.Range("E45").Value = '  -7.91%  '
$ws.Range("E46").Value = '  -3.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.98'
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("D48").Value = '2.707.42'
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.72'
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '66.96'
$ws.Range("E50").Value = '  -4.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.20'
$ws.Range("E51").Value = '  +1.97%  '
